# Insert a new data row at row 211 (pushing the existing rows 211-224 down to
# 212-225) and populate the new row with the latest weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 211-224 down to 212-225, inheriting formatting (style)
# from the row above, matching Excel's native "Insert" behaviour.
$ws.Rows.Item(211).Insert()

# Populate the newly inserted row 211 with the new observation.
$ws.Cells.Item(211, 1).Value  = 9
$ws.Cells.Item(211, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(211, 3).Value  = "Metropolitana"
$ws.Cells.Item(211, 4).Value  = 44578
$ws.Cells.Item(211, 5).Value  = 13
$ws.Cells.Item(211, 6).Value  = 100112021
$ws.Cells.Item(211, 7).Value  = "Ají"
$ws.Cells.Item(211, 8).Value  = "Americana (o)"
$ws.Cells.Item(211, 9).Value  = "Primera"
$ws.Cells.Item(211, 10).Value = 16
$ws.Cells.Item(211, 11).Value = 30000
$ws.Cells.Item(211, 12).Value = 32000
$ws.Cells.Item(211, 13).Value = 31000
$ws.Cells.Item(211, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(211, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(211, 16).Value = 1240
$ws.Cells.Item(211, 17).Value = 25
$ws.Cells.Item(211, 18).Value = "Hortaliza"
